## RBA 2.3 - Relatorio e Email
## Fix the placeholder text "QWREW"/"REW"/"Rew"/"rew" -> "QWR"/"QWER"/"Qwer"/"qwer"
## in the body ("A QWREW,") and in the page header (address block).

$d = $word.ActiveDocument

# wdReplaceOne = 1  (replace exactly the single match that Find just located)
function Replace-OneOccurrence {
    param($Range, $FindText, $ReplaceText, $MatchWholeWord)
    $r = $Range.Duplicate
    $ok = $r.Find.Execute($FindText, $true, $false, $MatchWholeWord, $false, $false, $true, 1, $false, $ReplaceText, 1)
    return $ok
}

# --- Main body: "A QWREW," -> "A QWR," -----------------------------------
$ok = Replace-OneOccurrence $d.Content "QWREW" "QWR" $false
Write-Host "body QWREW->QWR:" $ok

# --- Page header --------------------------------------------------------
$header = $d.Sections(1).Headers(1)

# "DIRETORIA DE ENSINO REGIAO REW" -> "... QWER"  (MatchWholeWord so the
# substring "REW" at the tail of "QWREW" a few characters later is skipped)
$ok = Replace-OneOccurrence $header.Range "REW" "QWER" $true
Write-Host "header REW->QWER:" $ok

# "QWREW - DEP." -> "QWR - DEP."
$ok = Replace-OneOccurrence $header.Range "QWREW" "QWR" $false
Write-Host "header QWREW->QWR:" $ok

# "Rew, no Rew - Rew - Rew - Rew" -> "Qwer, no Qwer - Qwer - Qwer - Qwer"
for ($i = 1; $i -le 5; $i++) {
    $ok = Replace-OneOccurrence $header.Range "Rew" "Qwer" $true
    Write-Host "header Rew->Qwer ($i):" $ok
}

# CEP / Tel / Email values: "rew" -> "qwer" (three occurrences)
for ($i = 1; $i -le 3; $i++) {
    $ok = Replace-OneOccurrence $header.Range "rew" "qwer" $true
    Write-Host "header rew->qwer ($i):" $ok
}
